$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows (10-13), which are no longer present in the updated data
$ws.Rows("10:13").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.113088
$ws.Range("H2").Value = 0.339264
$ws.Range("I2").Value = 0.3269336956678857
$ws.Range("J2").Value = 0.3269336956678857
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04067533333333333
$ws.Range("N2").Value = 0.122026
$ws.Range("O2").Value = 0.6477995434517174
$ws.Range("P2").Value = 0.6477995434517174
$ws.Range("Q2").Value = 0.004599892096
$ws.Range("R2").Value = 0.041399028864
$ws.Range("S2").Value = 0.2117874987926391
$ws.Range("T2").Value = 0.2117874987926391

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.113088
$ws.Range("H3").Value = 0.339264
$ws.Range("I3").Value = 0.3269336956678857
$ws.Range("J3").Value = 0.3269336956678857
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02211466666666667
$ws.Range("N3").Value = 0.066344
$ws.Range("O3").Value = 0.3522004565482826
$ws.Range("P3").Value = 0.3522004565482827
$ws.Range("Q3").Value = 0.002500903424
$ws.Range("R3").Value = 0.022508130816
$ws.Range("S3").Value = 0.1151461968752466
$ws.Range("T3").Value = 0.1151461968752467

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.084206
$ws.Range("H4").Value = 0.252618
$ws.Range("I4").Value = 0.2434367817753429
$ws.Range("J4").Value = 0.243436781775343
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04067533333333333
$ws.Range("N4").Value = 0.122026
$ws.Range("O4").Value = 0.6477995434517174
$ws.Range("P4").Value = 0.6477995434517174
$ws.Range("Q4").Value = 0.003425107118666667
$ws.Range("R4").Value = 0.030825964068
$ws.Range("S4").Value = 0.1576982360934225
$ws.Range("T4").Value = 0.1576982360934225

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.084206
$ws.Range("H5").Value = 0.252618
$ws.Range("I5").Value = 0.2434367817753429
$ws.Range("J5").Value = 0.243436781775343
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02211466666666667
$ws.Range("N5").Value = 0.066344
$ws.Range("O5").Value = 0.3522004565482826
$ws.Range("P5").Value = 0.3522004565482827
$ws.Range("Q5").Value = 0.001862187621333333
$ws.Range("R5").Value = 0.016759688592
$ws.Range("S5").Value = 0.08573854568192044
$ws.Range("T5").Value = 0.08573854568192046

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08859233333333333
$ws.Range("H6").Value = 0.265777
$ws.Range("I6").Value = 0.2561175274521424
$ws.Range("J6").Value = 0.2561175274521424
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04067533333333333
$ws.Range("N6").Value = 0.122026
$ws.Range("O6").Value = 0.6477995434517174
$ws.Range("P6").Value = 0.6477995434517174
$ws.Range("Q6").Value = 0.003603522689111111
$ws.Range("R6").Value = 0.032431704202
$ws.Range("S6").Value = 0.1659128173534806
$ws.Range("T6").Value = 0.1659128173534806

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08859233333333333
$ws.Range("H7").Value = 0.265777
$ws.Range("I7").Value = 0.2561175274521424
$ws.Range("J7").Value = 0.2561175274521424
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02211466666666667
$ws.Range("N7").Value = 0.066344
$ws.Range("O7").Value = 0.3522004565482826
$ws.Range("P7").Value = 0.3522004565482827
$ws.Range("Q7").Value = 0.001959189920888889
$ws.Range("R7").Value = 0.017632709288
$ws.Range("S7").Value = 0.09020471009866188
$ws.Range("T7").Value = 0.09020471009866189

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06001866666666666
$ws.Range("H8").Value = 0.180056
$ws.Range("I8").Value = 0.1735119951046289
$ws.Range("J8").Value = 0.1735119951046289
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04067533333333333
$ws.Range("N8").Value = 0.122026
$ws.Range("O8").Value = 0.6477995434517174
$ws.Range("P8").Value = 0.6477995434517174
$ws.Range("Q8").Value = 0.002441279272888889
$ws.Range("R8").Value = 0.021971513456
$ws.Range("S8").Value = 0.1124009912121752
$ws.Range("T8").Value = 0.1124009912121752

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06001866666666666
$ws.Range("H9").Value = 0.180056
$ws.Range("I9").Value = 0.1735119951046289
$ws.Range("J9").Value = 0.1735119951046289
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.02211466666666667
$ws.Range("N9").Value = 0.066344
$ws.Range("O9").Value = 0.3522004565482826
$ws.Range("P9").Value = 0.3522004565482827
$ws.Range("Q9").Value = 0.001327292807111111
$ws.Range("R9").Value = 0.011945635264
$ws.Range("S9").Value = 0.06111100389245368
$ws.Range("T9").Value = 0.0611110038924537
